# Updates the cryptocurrency price/volume table (rows 2-50) on Sheet1
# to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline updates: Price (col D) cells are plain text in the
# source data (values like "28.439.85" or "230.70"), so each is forced to
# Text format before assignment to stop Excel from re-interpreting it as a
# number, then the cell style is reset to "Normal" so no stray number format
# is left behind on the cell.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.439.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.55%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.787.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5867"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.86%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2753"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06698"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07542"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.793.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.04%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.760"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6052"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.30%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.030.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "75.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008657"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -10.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.430.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.392"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "206.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.52%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.762"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.47%  "

# Row 24
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.073"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1252"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.406"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06110"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.15%  "

# Row 31
$ws.Range("E31").Value = "  -1.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.758"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.761"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.671"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.84%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.041"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.73%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6390"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.503"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.42%  "

# Row 38
$ws.Range("E38").Value = "  -1.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.145.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.07%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01673"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.84%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.279"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8731"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.35%  "

# Row 43
$ws.Range("E43").Value = "  +0.32%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.938.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.59%  "

# Row 47
$ws.Range("E47").Value = "  -3.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.378"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.565"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05423"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.26%  "
